{"js": "// Insert a \"2008-07 \" date prefix into the document title so the\n// sub-project's month is mentioned in the heading, right before\n// \"System Object...\".\n// Net visible text change in the first (Heading 2) paragraph:\n//   \"Circle Language Spec Plan,<br>System Objects Spec,<br>Project Summary\"\n// becomes:\n//   \"Circle Language Spec Plan,<br>2008-07 System Objects Spec,<br>Project Summary\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The title is the first paragraph, styled \"Heading 2\".\nconst heading = paragraphs.items[0];\nheading.load(\"style\");\nawait context.sync();\n\nif (heading.style !== \"Heading 2\") {\n  throw new Error(\"Expected the first paragraph to use the Heading 2 style.\");\n}\n\n// Find \"System Object\" within that paragraph only (it also occurs later in\n// the document body) and insert the new date text immediately before it.\nconst matches = heading.search(\"System Object\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Could not find \"System Object\" in the title paragraph.');\n}\n\nmatches.items[0].insertText(\"2008-07 \", Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Insert the \"2008-07 \" month/date prefix into the document title so the\n# sub-project's month is mentioned in the heading, right before\n# \"System Object...\".\n#\n# Net visible text change in the first (Heading 2) paragraph:\n#   \"Circle Language Spec Plan,<br>System Objects Spec,<br>Project Summary\"\n# becomes:\n#   \"Circle Language Spec Plan,<br>2008-07 System Objects Spec,<br>Project Summary\"\n\n$d = $word.ActiveDocument\n\n$heading = $d.Paragraphs.Item(1)\nif ($heading.Style.NameLocal -ne \"Heading 2\") {\n    throw \"Expected the first paragraph to use the Heading 2 style.\"\n}\n\n# Search only within the title paragraph, since \"System Object\" also occurs\n# later in the document body.\n$matchRange = $heading.Range.Duplicate\n$find = $matchRange.Find\n$find.ClearFormatting()\n$find.Text = \"System Object\"\n$find.MatchCase = $true\n$find.Forward = $true\n$find.Wrap = 0\n\nif ($find.Execute()) {\n    $insertPoint = $d.Range($matchRange.Start, $matchRange.Start)\n    $insertPoint.InsertBefore(\"2008-07 \")\n} else {\n    throw 'Could not find \"System Object\" in the title paragraph.'\n}\n"}
